# Modified DSL for Tabbar, notification and adding new screenshots
#
# The "Stop Listening to the background intents" test case (row 10 of the
# TestCases sheet) is updated so its Steps/Description DSL now drives BOTH
# VT200-0438 and VT200-0439 (instead of VT200-0439 alone with a screenshot
# step), and finishes by returning Home, relaunching the app and checking
# for the "stoplistening" UI text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$nl = [char]10

# New "Description" (column H) DSL for row 10.
$h10Lines = @(
    'validate1',
    '{',
    'validate_PageTitle=Compliance JS specs',
    '};',
    'validate2',
    '{',
    'validate_PageTitle=Intent JS Test',
    '};',
    'validate3',
    '{',
    'validate_OldText_Exists=VT200-0438',
    '};',
    'validate4',
    '{',
    'validate_OldText_Exists=VT200-0439',
    '};',
    ''
)
$h10 = [string]::Join($nl, $h10Lines)

# New "Steps" (column G) DSL for row 10.
$g10Lines = @(
    'wait(3);',
    'validate1;',
    'link_Click(intent_test_link);',
    'validate2;',
    'SelectTestToRun(VT200_0438_string);',
    'ClickRunTest(runtest_top_xpath);',
    'validate3;',
    'ClickRunTest(runtest_bottom_xpath);',
    'wait(2);',
    'SelectTestToRun(VT200_0439_string);',
    'ClickRunTest(runtest_top_xpath);',
    'validate4;',
    'ClickRunTest(runtest_bottom_xpath);',
    'wait(2);',
    'press_Key(Home);',
    'launch_App_Device(com.rhomobile.compliancetest_js/com.rhomobile.rhodes.RhodesActivity);',
    'CheckUITextContains(stoplistening);'
)
$g10 = [string]::Join($nl, $g10Lines)

$ws.Range("H10").Value = $h10
$ws.Range("G10").Value = $g10

# Row grew taller to fit the extra validation/step lines.
$ws.Rows.Item(10).RowHeight = 243

# The saved view selection moved from D2 to D1.
$ws.Range("D1").Select() | Out-Null
